$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix header row formatting: drop the explicit row-level/yellow-ish
# fill style (fontId=1/fillId=4) that the header previously used, and
# restore the plain "fontId=1/no fill" style that the rest of the data
# cells use (style index 1 in styles.xml).
$ws.Rows.Item(1).ClearFormats()
$ws.Range("E2").Copy()
$ws.Range("A1:I1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Append the new daily records (16-21 Dec 2020) to the
# "Condicion_Pacientes" table, which grows the table/sheet from row 263
# to row 269.
$tbl = $ws.ListObjects.Item("Condicion_Pacientes")

$newData = @(
    @(16, 12, 2020, 4778, 607),
    @(17, 12, 2020, 4605, 627),
    @(18, 12, 2020, 4581, 533),
    @(19, 12, 2020, 1543, 170),
    @(20, 12, 2020, 667, 102),
    @(21, 12, 2020, 4669, 734)
)

foreach ($rowVals in $newData) {
    $lastDataRow = $tbl.Range.Row + $tbl.Range.Rows.Count - 1

    $listRow = $tbl.ListRows.Add()
    $r = $listRow.Range.Row

    # Copy the formatting (styles) of the previous data row down onto the
    # freshly-added row so it keeps the same look (fonts/number formats).
    $ws.Range("A" + $lastDataRow + ":I" + $lastDataRow).Copy()
    $ws.Range("A" + $r + ":I" + $r).PasteSpecial(-4122)
    $excel.CutCopyMode = $false
    $ws.Rows.Item($r).RowHeight = 14.25

    $ws.Cells.Item($r, 2).Value = $rowVals[0]
    $ws.Cells.Item($r, 3).Value = $rowVals[1]
    $ws.Cells.Item($r, 4).Value = $rowVals[2]
    $ws.Cells.Item($r, 1).Formula = "=+Hoja1!`$B$r&""/""&Hoja1!`$C$r&""/""&Hoja1!`$D$r"
    $ws.Cells.Item($r, 5).Value = $rowVals[3]
    $ws.Cells.Item($r, 6).Value = $rowVals[4]
}
